# Chart To Do List.xlsx - Data Statistics & Plot Work
# Calculated the adjusted carbon offset from the total carbon offset.
# Started plotting adjusted carbon offset by state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 ("Total Adjusted Carbon Offset by State") is now being worked on:
# mark the "Completed?" cell (E6) as "In Progress", matching the same
# highlight style already used for "In Progress" items (copy format from I6).
$ws.Range("I6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E6").Value2 = "In Progress"

# The Y-Axis title has been decided for every chart row (I6:I14), so mark
# the "Completed?3" column as "Yes" for all of them, matching the same
# highlight style already used for completed items (copy format from I5).
$ws.Range("I5").Copy() | Out-Null
$ws.Range("I6:I14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("I6:I14").Value2 = "Yes"

$excel.CutCopyMode = $false

# Move the active selection to D19, reflecting where work continued.
$ws.Range("D19").Select() | Out-Null
